$d = $word.ActiveDocument
$TAB = [char]9

function Replace-Text($searchText, $replaceText) {
    $r = $d.Content
    $found = $r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $found) {
        throw ("Replace-Text: not found: " + $searchText)
    }
    return $r
}

function Find-Only($searchText) {
    $r = $d.Content
    $found = $r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw ("Find-Only: not found: " + $searchText)
    }
    return $r
}

# ------------------------------------------------------------------
# 1. Update the record timestamp:
#    "THU Oct 12 11:10:12 PDT 2017" -> "SAT Oct 21 10:50:48 PDT 2017"
# ------------------------------------------------------------------
$dateSearch = "THU Oct 12 11:10:12 PDT 2017"
$dateReplace = "SAT Oct 21 10:50:48 PDT 2017"
Replace-Text $dateSearch $dateReplace | Out-Null

# ------------------------------------------------------------------
# 2. Person Name for this record: "- KSK" -> "- KNR HANUMANTHA"
# ------------------------------------------------------------------
$personSearch = "Person Name" + $TAB + $TAB + $TAB + $TAB + "- KSK"
$personReplace = "Person Name" + $TAB + $TAB + $TAB + $TAB + "- KNR HANUMANTHA"
Replace-Text $personSearch $personReplace | Out-Null

# ------------------------------------------------------------------
# 3. First "Item Name" value for this record: "- CARROT2" -> "- POTATO EVE 2"
# ------------------------------------------------------------------
$itemSearch = "Item Name" + $TAB + $TAB + $TAB + $TAB + "- CARROT2"
$itemReplace = "Item Name" + $TAB + $TAB + $TAB + $TAB + "- POTATO EVE 2"
$itemRange = Replace-Text $itemSearch $itemReplace
$itemNamePara = $itemRange.Paragraphs(1)
$afterIdx = $itemNamePara.Index

# ------------------------------------------------------------------
# 4. Insert the new "Amount Received" / "Amount Received mode" / blank /
#    "Item Name" paragraphs right after the paragraph we just edited.
# ------------------------------------------------------------------
$amtRecvLine = "Amount Received" + $TAB + $TAB + $TAB + "- 700"
$amtRecvModeLine = "Amount Received mode" + $TAB + $TAB + "- CASH AND CLEARD"
$itemLine2 = "Item Name" + $TAB + $TAB + $TAB + $TAB + "- POTATO EVE 2"

$itemNamePara.Range.InsertParagraphAfter()
$p1 = $d.Paragraphs($afterIdx + 1)
$p1.Range.Text = $amtRecvLine

$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs($afterIdx + 2)
$p2.Range.Text = $amtRecvModeLine

$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs($afterIdx + 3)
# p3 intentionally left blank

$p3.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs($afterIdx + 4)
$p4.Range.Text = $itemLine2

# ------------------------------------------------------------------
# 5. Number of KGs for this record: "- 67" -> "- 52"
# ------------------------------------------------------------------
$kgsSearch = "Number of KGs" + $TAB + $TAB + $TAB + "- 67"
$kgsReplace = "Number of KGs" + $TAB + $TAB + $TAB + "- 52"
Replace-Text $kgsSearch $kgsReplace | Out-Null

# ------------------------------------------------------------------
# 6. Rate for this record: "- 20" -> "- 13"
# ------------------------------------------------------------------
$rateSearch = "Rate" + $TAB + $TAB + $TAB + $TAB + $TAB + "- 20"
$rateReplace = "Rate" + $TAB + $TAB + $TAB + $TAB + $TAB + "- 13"
Replace-Text $rateSearch $rateReplace | Out-Null

# ------------------------------------------------------------------
# 7. Remove the "Transport & Miscellaneous" paragraph entirely for this record.
# ------------------------------------------------------------------
$transportSearch = "Transport & Miscellaneous" + $TAB + "- 60"
$transportRange = Find-Only $transportSearch
$transportPara = $transportRange.Paragraphs(1)
$transportPara.Range.Delete()

# ------------------------------------------------------------------
# 8. Total Price for this record: "- 1400.0" -> "- 676.0"
# ------------------------------------------------------------------
$totalSearch = "Total Price" + $TAB + $TAB + $TAB + $TAB + "- 1400.0"
$totalReplace = "Total Price" + $TAB + $TAB + $TAB + $TAB + "- 676.0"
$totalRange = Replace-Text $totalSearch $totalReplace
$totalPricePara = $totalRange.Paragraphs(1)
$totalIdx = $totalPricePara.Index

# ------------------------------------------------------------------
# 9. Amount balance for this record: un-bold it, update
#    "- 237385.0" -> "- 676.0", then add a new blank paragraph after it.
# ------------------------------------------------------------------
$balanceSearch = "Amount balance" + $TAB + $TAB + $TAB + "- 237385.0"
$balanceRange = Find-Only $balanceSearch
$balancePara = $balanceRange.Paragraphs(1)
$balancePara.Range.Delete()

$balanceLine = "Amount balance" + $TAB + $TAB + $TAB + "- 676.0"

$totalPricePara.Range.InsertParagraphAfter()
$pb1 = $d.Paragraphs($totalIdx + 1)
$pb1.Range.Text = $balanceLine

$pb1.Range.InsertParagraphAfter()
$pb2 = $d.Paragraphs($totalIdx + 2)
# pb2 intentionally left blank
